$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1. The paragraph that used to start with two manual line breaks
#    before "Here is the link to my Brewery UI-Project:" should only
#    keep a single line break.
# ---------------------------------------------------------------------
$d.Content.Find.Execute(
    "^l^lHere is the link to my Brewery UI-Project:",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "^lHere is the link to my Brewery UI-Project:", 2) | Out-Null

# ---------------------------------------------------------------------
# 2. Insert a brand new paragraph right after that one (and before the
#    existing BrewerySite.git hyperlink paragraph) containing the new
#    backup repository link, a manual line break, and the
#    "Here is the backup link:" label - all as plain text in one run.
# ---------------------------------------------------------------------
$labelPara = $null
$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    $candidate = $d.Paragraphs.Item($i)
    if ($candidate.Range.Text.Contains("Here is the link to my Brewery UI-Project:")) {
        $labelPara = $candidate
        break
    }
}

$labelRange = $labelPara.Range
$labelRange.InsertParagraphAfter()

$newPara = $labelPara.Next()
$newRange = $newPara.Range
$newRange.InsertBefore("https://github.com/Vexmage/Brewery-MSI.git")

$newRange2 = $newPara.Range
$lineBreak = [char]11
$newRange2.InsertAfter($lineBreak + "Here is the backup link:")

# ---------------------------------------------------------------------
# 3. Register a new "FollowedHyperlink" character style
#    (VisitedInternetLink) right next to the existing "Hyperlink"
#    character style (InternetLink), coloured maroon with a single
#    underline - mirroring the built-in FollowedHyperlink look.
# ---------------------------------------------------------------------
$styles = $d.Styles
$visited = $styles.Add("VisitedInternetLink", 2)
$visited.NameLocal = "FollowedHyperlink"
$visited.Font.Color = 128
$visited.Font.Underline = 1
